$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt2"
$ws.Range("C2").Value = "Fzd8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01070233333333333
$ws.Range("H2").Value = 0.032107
$ws.Range("I2").Value = 0.006017198313602724
$ws.Range("J2").Value = 0.006017198313602724
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.416422666666667
$ws.Range("N2").Value = 7.249268
$ws.Range("O2").Value = 0.2729115228630338
$ws.Range("P2").Value = 0.2729115228630338
$ws.Range("Q2").Value = 0.02586136085288889
$ws.Range("R2").Value = 0.232752247676
$ws.Range("S2").Value = 0.001642162755134198
$ws.Range("T2").Value = 0.001642162755134198

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt2"
$ws.Range("C3").Value = "Fzd8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01070233333333333
$ws.Range("H3").Value = 0.032107
$ws.Range("I3").Value = 0.006017198313602724
$ws.Range("J3").Value = 0.006017198313602724
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.913654666666667
$ws.Range("N3").Value = 8.740964
$ws.Range("O3").Value = 0.3290690586319826
$ws.Range("P3").Value = 0.3290690586319827
$ws.Range("Q3").Value = 0.03118290346088888
$ws.Range("R3").Value = 0.280646131148
$ws.Range("S3").Value = 0.001980073784659201
$ws.Range("T3").Value = 0.001980073784659202

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01070233333333333
$ws.Range("H4").Value = 0.032107
$ws.Range("I4").Value = 0.006017198313602724
$ws.Range("J4").Value = 0.006017198313602724
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.02921733333333333
$ws.Range("N4").Value = 0.08765200000000001
$ws.Range("O4").Value = 0.003299814657423431
$ws.Range("P4").Value = 0.003299814657423432
$ws.Range("Q4").Value = 0.0003126936404444444
$ws.Range("R4").Value = 0.002814242764
$ws.Range("S4").Value = 0.00001985563919184982
$ws.Range("T4").Value = 0.00001985563919184982

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd8"
$ws.Range("D5").Value = "Neutro"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01070233333333333
$ws.Range("H5").Value = 0.032107
$ws.Range("I5").Value = 0.006017198313602724
$ws.Range("J5").Value = 0.006017198313602724
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02851766666666667
$ws.Range("N5").Value = 0.085553
$ws.Range("O5").Value = 0.003220794087830817
$ws.Range("P5").Value = 0.003220794087830818
$ws.Range("Q5").Value = 0.0003052055745555555
$ws.Range("R5").Value = 0.002746850171
$ws.Range("S5").Value = 0.00001938015675375722
$ws.Range("T5").Value = 0.00001938015675375722

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd8"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01070233333333333
$ws.Range("H6").Value = 0.032107
$ws.Range("I6").Value = 0.006017198313602724
$ws.Range("J6").Value = 0.006017198313602724
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.466422333333333
$ws.Range("N6").Value = 10.399267
$ws.Range("O6").Value = 0.3914988097597293
$ws.Range("P6").Value = 0.3914988097597293
$ws.Range("Q6").Value = 0.03709880728544444
$ws.Range("R6").Value = 0.3338892655689999
$ws.Range("S6").Value = 0.002355725977863717
$ws.Range("T6").Value = 0.002355725977863717

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt2"
$ws.Range("C7").Value = "Fzd8"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.767921666666667
$ws.Range("H7").Value = 5.303765
$ws.Range("I7").Value = 0.9939828016863973
$ws.Range("J7").Value = 0.9939828016863973
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.416422666666667
$ws.Range("N7").Value = 7.249268
$ws.Range("O7").Value = 0.2729115228630338
$ws.Range("P7").Value = 0.2729115228630338
$ws.Range("Q7").Value = 4.272045988224445
$ws.Range("R7").Value = 38.44841389402
$ws.Range("S7").Value = 0.2712693601078996
$ws.Range("T7").Value = 0.2712693601078996

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt2"
$ws.Range("C8").Value = "Fzd8"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.767921666666667
$ws.Range("H8").Value = 5.303765
$ws.Range("I8").Value = 0.9939828016863973
$ws.Range("J8").Value = 0.9939828016863973
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.913654666666667
$ws.Range("N8").Value = 8.740964
$ws.Range("O8").Value = 0.3290690586319826
$ws.Range("P8").Value = 0.3290690586319827
$ws.Range("Q8").Value = 5.151113214384445
$ws.Range("R8").Value = 46.36001892946
$ws.Range("S8").Value = 0.3270889848473234
$ws.Range("T8").Value = 0.3270889848473235

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt2"
$ws.Range("C9").Value = "Fzd8"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.767921666666667
$ws.Range("H9").Value = 5.303765
$ws.Range("I9").Value = 0.9939828016863973
$ws.Range("J9").Value = 0.9939828016863973
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.02921733333333333
$ws.Range("N9").Value = 0.08765200000000001
$ws.Range("O9").Value = 0.003299814657423431
$ws.Range("P9").Value = 0.003299814657423432
$ws.Range("Q9").Value = 0.05165395664222223
$ws.Range("R9").Value = 0.46488560978
$ws.Range("S9").Value = 0.003279959018231581
$ws.Range("T9").Value = 0.003279959018231582

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt2"
$ws.Range("C10").Value = "Fzd8"
$ws.Range("D10").Value = "Neutro"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.767921666666667
$ws.Range("H10").Value = 5.303765
$ws.Range("I10").Value = 0.9939828016863973
$ws.Range("J10").Value = 0.9939828016863973
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02851766666666667
$ws.Range("N10").Value = 0.085553
$ws.Range("O10").Value = 0.003220794087830817
$ws.Range("P10").Value = 0.003220794087830818
$ws.Range("Q10").Value = 0.05041700078277778
$ws.Range("R10").Value = 0.453753007045
$ws.Range("S10").Value = 0.00320141393107706
$ws.Range("T10").Value = 0.003201413931077061

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt2"
$ws.Range("C11").Value = "Fzd8"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.767921666666667
$ws.Range("H11").Value = 5.303765
$ws.Range("I11").Value = 0.9939828016863973
$ws.Range("J11").Value = 0.9939828016863973
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.466422333333333
$ws.Range("N11").Value = 10.399267
$ws.Range("O11").Value = 0.3914988097597293
$ws.Range("P11").Value = 0.3914988097597293
$ws.Range("Q11").Value = 6.128363148917223
$ws.Range("R11").Value = 55.15526834025501
$ws.Range("S11").Value = 0.3891430837818656
$ws.Range("T11").Value = 0.3891430837818656

